# "New Branch for work" - add Sheet2 "ShopKeeperForm"-style task column (G) with new
# stories, resize a couple of rows to fit the new wrapped text, highlight Sheet1 G6 in
# purple, and leave the selections where the author left them.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet2: new "Tasks" entries in column G -------------------------------------
# Give G3:G8 the same wrap/centered style already used by column D (xf s="14"),
# then fill in the values. The values are written in the same order the author
# originally typed them in (this controls shared-string insertion order), not in
# top-to-bottom row order.
$ws2.Range("D3").Copy()
$ws2.Range("G3:G8").PasteSpecial(-4122)

$ws2.Range("G8").Value = "pushing to master"
$ws2.Range("G7").Value = "completing the first version "
$ws2.Range("G4").Value = "make Shopkeeper project"
$ws2.Range("G3").Value = "Testing github push and fetch"
$ws2.Range("G5").Value = "desing shopkeeper form"
$ws2.Range("G6").Value = "write code for applications buttons"

# Row 4 now needs two lines of wrapped text, and the new row 8 needs room too.
$ws2.Range("A4").RowHeight = 45
$ws2.Range("A8").RowHeight = 30

# --- Sheet1: highlight G6 with the purple fill/font combo ------------------------
$ws1.Range("G6").Font.Color = 10498160
$ws1.Range("G6").Interior.Color = 10498160

# --- Restore selections ------------------------------------------------------------
# Select Sheet2's cell first, then Sheet1's, so Sheet1 ends up the active tab again
# (matching the workbook's original tabSelected state).
[void]$ws2.Range("G7").Select()
[void]$ws1.Range("I16:J16").Select()

# --- Calculation options (best-effort; engine recomputes calcPr on save) ---------
$excel.Iteration = $true
$excel.MaxChange = 0.0001
$excel.MultiThreadedCalculation = 0
